# Update to query, include AggregationIntervalUnitCV
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Time Range": add AggregationIntervalUnitCV (col B) and
# numUniqueYears (col E) columns into the summary table, shifting the
# existing minYear / maxYear / TimeRange_Yrs data and the concatenated
# "pipe" helper formulas over to make room for them.
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Time Range")

# Header row values
$tr.Range("B1").Value = "AggregationIntervalUnitCV "
$tr.Range("C1").Value = "minYear"
$tr.Range("D1").Value = "maxYear"
$tr.Range("E1").Value = "numUniqueYears"
$tr.Range("F1").Value = "TimeRange_Yrs"

# F1 (TimeRange_Yrs) keeps the bold header formatting that used to live on D1;
# D1 reverts to the default (unbolded) style.
$tr.Range("D1").Copy() | Out-Null
$tr.Range("F1").PasteSpecial(-4122) | Out-Null
$tr.Range("D1").ClearFormats() | Out-Null
$tr.Range("F1").Value = "TimeRange_Yrs"

# Data rows: State | AggregationIntervalUnitCV | minYear | maxYear | numUniqueYears
$tr.Range("B2").Value = "Monthly"
$tr.Range("C2").Value = 2013
$tr.Range("D2").Value = 2016
$tr.Range("E2").Value = 4

$tr.Range("B3").Value = "Monthly"
$tr.Range("C3").Value = 1990
$tr.Range("D3").Value = 2018
$tr.Range("E3").Value = 29

$tr.Range("B4").Value = "Annual"
$tr.Range("C4").Value = 2010
$tr.Range("D4").Value = 2015
$tr.Range("E4").Value = 2

$tr.Range("B5").Value = "Monthly"
$tr.Range("C5").Value = 1955
$tr.Range("D5").Value = 2021
$tr.Range("E5").Value = 67

$tr.Range("B6").Value = "Year"
$tr.Range("C6").Value = 1957
$tr.Range("D6").Value = 2020
$tr.Range("E6").Value = 64

# TimeRange_Yrs (col F) = maxYear - minYear; F3:F6 fill down as one shared formula
$tr.Range("F2").Formula = "=D2-C2"
$tr.Range("F3:F6").Formula = "=D3-C3"

# Pipe-concatenated helper column (was H, now I)
$tr.Range("I1").Formula = '=A1&"|"&B1&"|"&C1&"|"&D1&"|"&E1&"|"&F1'

$tr.Range("H2").Copy() | Out-Null
$tr.Range("I2").PasteSpecial(-4122) | Out-Null
$tr.Range("I2").Value = "----|----|---- |----|----|----"

$tr.Range("I3").Formula = '=A2&"|"&B2&"|"&C2&"|"&D2&"|"&E2&"|"&F2'
$tr.Range("I4:I7").Formula = '=A3&"|"&B3&"|"&C3&"|"&D3&"|"&E3&"|"&F3'

# Leftover bold-format artifact on F9 (stray paste/drag from F1), mirroring
# the same kind of blank styled cell left on "POD v POU"!F16.
$tr.Range("F1").Copy() | Out-Null
$tr.Range("F9").PasteSpecial(-4122) | Out-Null
$tr.Range("F9").ClearContents() | Out-Null

# The old helper column H is now superseded by column I; clear it entirely
# (values, formulas and column formatting).
$tr.Range("H1:H9").Clear() | Out-Null

$tr.PageSetup.Orientation = 1

# Move the active selection (purely a view-state change in the source file)
$tr.Range("K18").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "POD v POU": data is unchanged, but the per-row H formulas collapse
# into a single shared formula group (H3:H8) in the saved file; replicate by
# writing the formula across the whole range at once.
# ---------------------------------------------------------------------------
$pod = $wb.Worksheets.Item("POD v POU")
$pod.Range("H3:H8").Formula = '=A2&"|"&B2&"|"&C2'

$pod.Range("J26:J27").Select() | Out-Null
